$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): add date / legislator_name / legislator_id columns ---
$ws.Cells.Item(1,8).Value = "date"
$ws.Cells.Item(1,9).Value = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"

# match formatting of the existing header cell (G1, style index 1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

# --- Data row (row 2): fill in the new values ---
# Force the date to be stored as text (matches "2013-12-03" literal string in source)
$ws.Range("H2").NumberFormat = "@"
$ws.Cells.Item(2,8).Value = "2013-12-03"
$ws.Cells.Item(2,9).Value = "邱議瑩"
$ws.Cells.Item(2,10).Value = 913

# match formatting of the existing data cell (G2, style index 2)
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2:J2").PasteSpecial(-4122) | Out-Null
